# --- Refresh the "panel_query_time" timestamps on the "data" sheet (F2:F10) ---
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:20:16.826464"
$data.Range("F3").Value = "2021-10-05 14:20:16.826472"
$data.Range("F4").Value = "2021-10-05 14:20:16.826475"
$data.Range("F5").Value = "2021-10-05 14:20:16.826478"
$data.Range("F6").Value = "2021-10-05 14:20:16.826481"
$data.Range("F7").Value = "2021-10-05 14:20:16.826483"
$data.Range("F8").Value = "2021-10-05 14:20:16.826486"
$data.Range("F9").Value = "2021-10-05 14:20:16.826489"
$data.Range("F10").Value = "2021-10-05 14:20:16.826491"

# --- Add the new "metadata" sheet, placed immediately after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Reuse the bold/centered/bordered header formatting from the "data" sheet
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Familial hypoparathyroidism"
$meta.Range("C2").Value = 312
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.5"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2020-12-07T14:31:40.994962Z"
$meta.Range("F2").Value = "2021-10-05 14:20:16.823063"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/312/?format=json"

$excel.CutCopyMode = $false
